$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 194
$ws.Range("I2").Value = 586
$ws.Range("J2").Value = 2279
$ws.Range("K2").Value = 18
$ws.Range("L2").Value = 639
$ws.Range("M2").Value = 27
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 33
$ws.Range("S2").Value = 257
$ws.Range("T2").Value = 406
$ws.Range("V2").Value = 3662
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3647
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 68
$ws.Range("AA2").Value = 21
